$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intra-Session Statistics")

# Update headers in both tables: D3/R3 first (so "Resp" becomes shared-string 37),
# then C3/Q3 ("Group" becomes shared-string 38)
$ws.Range("D3").Value = "Resp"
$ws.Range("R3").Value = "Resp"
$ws.Range("C3").Value = "Group"
$ws.Range("Q3").Value = "Group"

# Existing rows 4-23 become "Group 1" rows: C/Q = group number 1, D/R = old C/Q value (Resp)
$ws.Range("D4").Value = 4.5551000000000004
$ws.Range("C4").Value = 1
$ws.Range("C4").Style = "Normal"
$ws.Range("R4").Value = 4.1835000000000004
$ws.Range("Q4").Value = 1
$ws.Range("Q4").Style = "Normal"
$ws.Range("D5").Value = 4.5209999999999999
$ws.Range("C5").Value = 1
$ws.Range("C5").Style = "Normal"
$ws.Range("R5").Value = 4.1134000000000004
$ws.Range("Q5").Value = 1
$ws.Range("Q5").Style = "Normal"
$ws.Range("D6").Value = 4.8205
$ws.Range("C6").Value = 1
$ws.Range("C6").Style = "Normal"
$ws.Range("R6").Value = 4.1486000000000001
$ws.Range("Q6").Value = 1
$ws.Range("Q6").Style = "Normal"
$ws.Range("D7").Value = 4.9162999999999997
$ws.Range("C7").Value = 1
$ws.Range("C7").Style = "Normal"
$ws.Range("R7").Value = 3.1766000000000001
$ws.Range("Q7").Value = 1
$ws.Range("Q7").Style = "Normal"
$ws.Range("D8").Value = 4.8445
$ws.Range("C8").Value = 1
$ws.Range("C8").Style = "Normal"
$ws.Range("R8").Value = 4.1554000000000002
$ws.Range("Q8").Value = 1
$ws.Range("Q8").Style = "Normal"
$ws.Range("D9").Value = 2.1589
$ws.Range("C9").Value = 1
$ws.Range("C9").Style = "Normal"
$ws.Range("R9").Value = 3.3216000000000001
$ws.Range("Q9").Value = 1
$ws.Range("Q9").Style = "Normal"
$ws.Range("D10").Value = 2.6341000000000001
$ws.Range("C10").Value = 1
$ws.Range("C10").Style = "Normal"
$ws.Range("R10").Value = 3.3795999999999999
$ws.Range("Q10").Value = 1
$ws.Range("Q10").Style = "Normal"
$ws.Range("D11").Value = 2.5419999999999998
$ws.Range("C11").Value = 1
$ws.Range("C11").Style = "Normal"
$ws.Range("R11").Value = 3.302
$ws.Range("Q11").Value = 1
$ws.Range("Q11").Style = "Normal"
$ws.Range("D12").Value = 2.9980000000000002
$ws.Range("C12").Value = 1
$ws.Range("C12").Style = "Normal"
$ws.Range("R12").Value = 3.4218000000000002
$ws.Range("Q12").Value = 1
$ws.Range("Q12").Style = "Normal"
$ws.Range("D13").Value = 2.8119999999999998
$ws.Range("C13").Value = 1
$ws.Range("C13").Style = "Normal"
$ws.Range("R13").Value = 3.605
$ws.Range("Q13").Value = 1
$ws.Range("Q13").Style = "Normal"
$ws.Range("D14").Value = 2.8096999999999999
$ws.Range("C14").Value = 1
$ws.Range("C14").Style = "Normal"
$ws.Range("R14").Value = 3.2883
$ws.Range("Q14").Value = 1
$ws.Range("Q14").Style = "Normal"
$ws.Range("D15").Value = 3.1785000000000001
$ws.Range("C15").Value = 1
$ws.Range("C15").Style = "Normal"
$ws.Range("R15").Value = 3.2995000000000001
$ws.Range("Q15").Value = 1
$ws.Range("Q15").Style = "Normal"
$ws.Range("D16").Value = 3.4561000000000002
$ws.Range("C16").Value = 1
$ws.Range("C16").Style = "Normal"
$ws.Range("R16").Value = 3.2530000000000001
$ws.Range("Q16").Value = 1
$ws.Range("Q16").Style = "Normal"
$ws.Range("D17").Value = 3.7456
$ws.Range("C17").Value = 1
$ws.Range("C17").Style = "Normal"
$ws.Range("R17").Value = 3.1150000000000002
$ws.Range("Q17").Value = 1
$ws.Range("Q17").Style = "Normal"
$ws.Range("D18").Value = 2.6009000000000002
$ws.Range("C18").Value = 1
$ws.Range("C18").Style = "Normal"
$ws.Range("R18").Value = 3.1488999999999998
$ws.Range("Q18").Value = 1
$ws.Range("Q18").Style = "Normal"
$ws.Range("D19").Value = 1.6719999999999999
$ws.Range("C19").Value = 1
$ws.Range("C19").Style = "Normal"
$ws.Range("R19").Value = 2.6928000000000001
$ws.Range("Q19").Value = 1
$ws.Range("Q19").Style = "Normal"
$ws.Range("D20").Value = 2.2176999999999998
$ws.Range("C20").Value = 1
$ws.Range("C20").Style = "Normal"
$ws.Range("R20").Value = 3.0409000000000002
$ws.Range("Q20").Value = 1
$ws.Range("Q20").Style = "Normal"
$ws.Range("D21").Value = 2.4321000000000002
$ws.Range("C21").Value = 1
$ws.Range("C21").Style = "Normal"
$ws.Range("R21").Value = 2.5383
$ws.Range("Q21").Value = 1
$ws.Range("Q21").Style = "Normal"
$ws.Range("D22").Value = 2.7351000000000001
$ws.Range("C22").Value = 1
$ws.Range("C22").Style = "Normal"
$ws.Range("R22").Value = 2.5951
$ws.Range("Q22").Value = 1
$ws.Range("Q22").Style = "Normal"
$ws.Range("D23").Value = 2.1021000000000001
$ws.Range("C23").Value = 1
$ws.Range("C23").Style = "Normal"
$ws.Range("R23").Value = 2.6928000000000001
$ws.Range("Q23").Value = 1
$ws.Range("Q23").Style = "Normal"

# New rows 24-43: Group 2 data (old D/R values moved to new rows)
# A/B/D and O/P/R keep the centered style (xlCenter) used by the rest of the table;
# C/Q (group number) stay with default/no style, matching the existing Group-number cells.
$ws.Range("A24").Value = 11002
$ws.Range("A24").HorizontalAlignment = -4108
$ws.Range("B24").Value = 1
$ws.Range("B24").HorizontalAlignment = -4108
$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 3.984
$ws.Range("D24").HorizontalAlignment = -4108
$ws.Range("O24").Value = 11002
$ws.Range("O24").HorizontalAlignment = -4108
$ws.Range("P24").Value = 1
$ws.Range("P24").HorizontalAlignment = -4108
$ws.Range("Q24").Value = 2
$ws.Range("R24").Value = 4.5785
$ws.Range("R24").HorizontalAlignment = -4108
$ws.Range("A25").Value = 11002
$ws.Range("A25").HorizontalAlignment = -4108
$ws.Range("B25").Value = 2
$ws.Range("B25").HorizontalAlignment = -4108
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 4.2622999999999998
$ws.Range("D25").HorizontalAlignment = -4108
$ws.Range("O25").Value = 11002
$ws.Range("O25").HorizontalAlignment = -4108
$ws.Range("P25").Value = 2
$ws.Range("P25").HorizontalAlignment = -4108
$ws.Range("Q25").Value = 2
$ws.Range("R25").Value = 4.5646000000000004
$ws.Range("R25").HorizontalAlignment = -4108
$ws.Range("A26").Value = 11002
$ws.Range("A26").HorizontalAlignment = -4108
$ws.Range("B26").Value = 3
$ws.Range("B26").HorizontalAlignment = -4108
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 4.6292999999999997
$ws.Range("D26").HorizontalAlignment = -4108
$ws.Range("O26").Value = 11002
$ws.Range("O26").HorizontalAlignment = -4108
$ws.Range("P26").Value = 3
$ws.Range("P26").HorizontalAlignment = -4108
$ws.Range("Q26").Value = 2
$ws.Range("R26").Value = 4.1680999999999999
$ws.Range("R26").HorizontalAlignment = -4108
$ws.Range("A27").Value = 11002
$ws.Range("A27").HorizontalAlignment = -4108
$ws.Range("B27").Value = 4
$ws.Range("B27").HorizontalAlignment = -4108
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 4.9283000000000001
$ws.Range("D27").HorizontalAlignment = -4108
$ws.Range("O27").Value = 11002
$ws.Range("O27").HorizontalAlignment = -4108
$ws.Range("P27").Value = 4
$ws.Range("P27").HorizontalAlignment = -4108
$ws.Range("Q27").Value = 2
$ws.Range("R27").Value = 4.0476999999999999
$ws.Range("R27").HorizontalAlignment = -4108
$ws.Range("A28").Value = 11002
$ws.Range("A28").HorizontalAlignment = -4108
$ws.Range("B28").Value = 5
$ws.Range("B28").HorizontalAlignment = -4108
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 5.0427999999999997
$ws.Range("D28").HorizontalAlignment = -4108
$ws.Range("O28").Value = 11002
$ws.Range("O28").HorizontalAlignment = -4108
$ws.Range("P28").Value = 5
$ws.Range("P28").HorizontalAlignment = -4108
$ws.Range("Q28").Value = 2
$ws.Range("R28").Value = 4.0164
$ws.Range("R28").HorizontalAlignment = -4108
$ws.Range("A29").Value = 11057
$ws.Range("A29").HorizontalAlignment = -4108
$ws.Range("B29").Value = 1
$ws.Range("B29").HorizontalAlignment = -4108
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 2.972
$ws.Range("D29").HorizontalAlignment = -4108
$ws.Range("O29").Value = 11057
$ws.Range("O29").HorizontalAlignment = -4108
$ws.Range("P29").Value = 1
$ws.Range("P29").HorizontalAlignment = -4108
$ws.Range("Q29").Value = 2
$ws.Range("R29").Value = 3.5661999999999998
$ws.Range("R29").HorizontalAlignment = -4108
$ws.Range("A30").Value = 11057
$ws.Range("A30").HorizontalAlignment = -4108
$ws.Range("B30").Value = 2
$ws.Range("B30").HorizontalAlignment = -4108
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 2.3220000000000001
$ws.Range("D30").HorizontalAlignment = -4108
$ws.Range("O30").Value = 11057
$ws.Range("O30").HorizontalAlignment = -4108
$ws.Range("P30").Value = 2
$ws.Range("P30").HorizontalAlignment = -4108
$ws.Range("Q30").Value = 2
$ws.Range("R30").Value = 3.59
$ws.Range("R30").HorizontalAlignment = -4108
$ws.Range("A31").Value = 11057
$ws.Range("A31").HorizontalAlignment = -4108
$ws.Range("B31").Value = 3
$ws.Range("B31").HorizontalAlignment = -4108
$ws.Range("C31").Value = 2
$ws.Range("D31").Value = 2.7309000000000001
$ws.Range("D31").HorizontalAlignment = -4108
$ws.Range("O31").Value = 11057
$ws.Range("O31").HorizontalAlignment = -4108
$ws.Range("P31").Value = 3
$ws.Range("P31").HorizontalAlignment = -4108
$ws.Range("Q31").Value = 2
$ws.Range("R31").Value = 3.1556999999999999
$ws.Range("R31").HorizontalAlignment = -4108
$ws.Range("A32").Value = 11057
$ws.Range("A32").HorizontalAlignment = -4108
$ws.Range("B32").Value = 4
$ws.Range("B32").HorizontalAlignment = -4108
$ws.Range("C32").Value = 2
$ws.Range("D32").Value = 2.8187000000000002
$ws.Range("D32").HorizontalAlignment = -4108
$ws.Range("O32").Value = 11057
$ws.Range("O32").HorizontalAlignment = -4108
$ws.Range("P32").Value = 4
$ws.Range("P32").HorizontalAlignment = -4108
$ws.Range("Q32").Value = 2
$ws.Range("R32").Value = 3.1133999999999999
$ws.Range("R32").HorizontalAlignment = -4108
$ws.Range("A33").Value = 11057
$ws.Range("A33").HorizontalAlignment = -4108
$ws.Range("B33").Value = 5
$ws.Range("B33").HorizontalAlignment = -4108
$ws.Range("C33").Value = 2
$ws.Range("D33").Value = 2.9215
$ws.Range("D33").HorizontalAlignment = -4108
$ws.Range("O33").Value = 11057
$ws.Range("O33").HorizontalAlignment = -4108
$ws.Range("P33").Value = 5
$ws.Range("P33").HorizontalAlignment = -4108
$ws.Range("Q33").Value = 2
$ws.Range("R33").Value = 3.58
$ws.Range("R33").HorizontalAlignment = -4108
$ws.Range("A34").Value = 11108
$ws.Range("A34").HorizontalAlignment = -4108
$ws.Range("B34").Value = 1
$ws.Range("B34").HorizontalAlignment = -4108
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 2.7021000000000002
$ws.Range("D34").HorizontalAlignment = -4108
$ws.Range("O34").Value = 11108
$ws.Range("O34").HorizontalAlignment = -4108
$ws.Range("P34").Value = 1
$ws.Range("P34").HorizontalAlignment = -4108
$ws.Range("Q34").Value = 2
$ws.Range("R34").Value = 2.911
$ws.Range("R34").HorizontalAlignment = -4108
$ws.Range("A35").Value = 11108
$ws.Range("A35").HorizontalAlignment = -4108
$ws.Range("B35").Value = 2
$ws.Range("B35").HorizontalAlignment = -4108
$ws.Range("C35").Value = 2
$ws.Range("D35").Value = 3.2
$ws.Range("D35").HorizontalAlignment = -4108
$ws.Range("O35").Value = 11108
$ws.Range("O35").HorizontalAlignment = -4108
$ws.Range("P35").Value = 2
$ws.Range("P35").HorizontalAlignment = -4108
$ws.Range("Q35").Value = 2
$ws.Range("R35").Value = 2.8378999999999999
$ws.Range("R35").HorizontalAlignment = -4108
$ws.Range("A36").Value = 11108
$ws.Range("A36").HorizontalAlignment = -4108
$ws.Range("B36").Value = 3
$ws.Range("B36").HorizontalAlignment = -4108
$ws.Range("C36").Value = 2
$ws.Range("D36").Value = 3.2254999999999998
$ws.Range("D36").HorizontalAlignment = -4108
$ws.Range("O36").Value = 11108
$ws.Range("O36").HorizontalAlignment = -4108
$ws.Range("P36").Value = 3
$ws.Range("P36").HorizontalAlignment = -4108
$ws.Range("Q36").Value = 2
$ws.Range("R36").Value = 3.1657000000000002
$ws.Range("R36").HorizontalAlignment = -4108
$ws.Range("A37").Value = 11108
$ws.Range("A37").HorizontalAlignment = -4108
$ws.Range("B37").Value = 4
$ws.Range("B37").HorizontalAlignment = -4108
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = 2.9350000000000001
$ws.Range("D37").HorizontalAlignment = -4108
$ws.Range("O37").Value = 11108
$ws.Range("O37").HorizontalAlignment = -4108
$ws.Range("P37").Value = 4
$ws.Range("P37").HorizontalAlignment = -4108
$ws.Range("Q37").Value = 2
$ws.Range("R37").Value = 3.3003
$ws.Range("R37").HorizontalAlignment = -4108
$ws.Range("A38").Value = 11108
$ws.Range("A38").HorizontalAlignment = -4108
$ws.Range("B38").Value = 5
$ws.Range("B38").HorizontalAlignment = -4108
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 3.1173999999999999
$ws.Range("D38").HorizontalAlignment = -4108
$ws.Range("O38").Value = 11108
$ws.Range("O38").HorizontalAlignment = -4108
$ws.Range("P38").Value = 5
$ws.Range("P38").HorizontalAlignment = -4108
$ws.Range("Q38").Value = 2
$ws.Range("R38").Value = 3.7035
$ws.Range("R38").HorizontalAlignment = -4108
$ws.Range("A39").Value = 11112
$ws.Range("A39").HorizontalAlignment = -4108
$ws.Range("B39").Value = 1
$ws.Range("B39").HorizontalAlignment = -4108
$ws.Range("C39").Value = 2
$ws.Range("D39").Value = 2.3041
$ws.Range("D39").HorizontalAlignment = -4108
$ws.Range("O39").Value = 11112
$ws.Range("O39").HorizontalAlignment = -4108
$ws.Range("P39").Value = 1
$ws.Range("P39").HorizontalAlignment = -4108
$ws.Range("Q39").Value = 2
$ws.Range("R39").Value = 2.7675999999999998
$ws.Range("R39").HorizontalAlignment = -4108
$ws.Range("A40").Value = 11112
$ws.Range("A40").HorizontalAlignment = -4108
$ws.Range("B40").Value = 2
$ws.Range("B40").HorizontalAlignment = -4108
$ws.Range("C40").Value = 2
$ws.Range("D40").Value = 2.0830000000000002
$ws.Range("D40").HorizontalAlignment = -4108
$ws.Range("O40").Value = 11112
$ws.Range("O40").HorizontalAlignment = -4108
$ws.Range("P40").Value = 2
$ws.Range("P40").HorizontalAlignment = -4108
$ws.Range("Q40").Value = 2
$ws.Range("R40").Value = 2.2237
$ws.Range("R40").HorizontalAlignment = -4108
$ws.Range("A41").Value = 11112
$ws.Range("A41").HorizontalAlignment = -4108
$ws.Range("B41").Value = 3
$ws.Range("B41").HorizontalAlignment = -4108
$ws.Range("C41").Value = 2
$ws.Range("D41").Value = 2.3456999999999999
$ws.Range("D41").HorizontalAlignment = -4108
$ws.Range("O41").Value = 11112
$ws.Range("O41").HorizontalAlignment = -4108
$ws.Range("P41").Value = 3
$ws.Range("P41").HorizontalAlignment = -4108
$ws.Range("Q41").Value = 2
$ws.Range("R41").Value = 2.6067999999999998
$ws.Range("R41").HorizontalAlignment = -4108
$ws.Range("A42").Value = 11112
$ws.Range("A42").HorizontalAlignment = -4108
$ws.Range("B42").Value = 4
$ws.Range("B42").HorizontalAlignment = -4108
$ws.Range("C42").Value = 2
$ws.Range("D42").Value = 2.6530999999999998
$ws.Range("D42").HorizontalAlignment = -4108
$ws.Range("O42").Value = 11112
$ws.Range("O42").HorizontalAlignment = -4108
$ws.Range("P42").Value = 4
$ws.Range("P42").HorizontalAlignment = -4108
$ws.Range("Q42").Value = 2
$ws.Range("R42").Value = 3.1389999999999998
$ws.Range("R42").HorizontalAlignment = -4108
$ws.Range("A43").Value = 11112
$ws.Range("A43").HorizontalAlignment = -4108
$ws.Range("B43").Value = 5
$ws.Range("B43").HorizontalAlignment = -4108
$ws.Range("C43").Value = 2
$ws.Range("D43").Value = 2.9794999999999998
$ws.Range("D43").HorizontalAlignment = -4108
$ws.Range("O43").Value = 11112
$ws.Range("O43").HorizontalAlignment = -4108
$ws.Range("P43").Value = 5
$ws.Range("P43").HorizontalAlignment = -4108
$ws.Range("Q43").Value = 2
$ws.Range("R43").Value = 2.4180999999999999
$ws.Range("R43").HorizontalAlignment = -4108

# Selection + view changes
$ws.Range("K12").Select() | Out-Null

# Page setup (portrait orientation, matches the saved print settings)
$ws.PageSetup.Orientation = 1
